$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert previously-empty-string cells into numeric 1 values ---
$ws.Range("H4").Value = 1
$ws.Range("L4").Value = 1

$ws.Range("H5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("L5").Value = 1

$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("L7").Value = 1

$ws.Range("H9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("L9").Value = 1

$ws.Range("H12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("L12").Value = 1

# --- Add summary rows 18 and 19 with COUNTIF formulas ---
$ws.Range("H18").Formula = "=COUNTIF(H2:H17,1)"
$ws.Range("J18").Formula = "=COUNTIF(J2:J17,1)"
$ws.Range("L18").Formula = "=COUNTIF(L2:L17,1)"

$ws.Range("H19").Formula = "=COUNTIF(H2:H17,0)"
$ws.Range("J19").Formula = "=COUNTIF(J2:J17,0)"
$ws.Range("L19").Formula = "=COUNTIF(L2:L17,0)"

# --- Update sheet view state ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("I23").Select()
